$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (event/strategy fields regenerated)
$ws.Range("A2").Value = 0
$ws.Range("C2").Value = 0.8398746252059937
$ws.Range("B2").Value = "Remorse"
$ws.Range("D2").Value = "BecomeRich"
$ws.Range("E2").Value = "Situation Modification"

# Strategy related label changes
$ws.Range("G11").Value = "[Cognitive Change, Strongly]"

# Move the dominant-personality value down one row, freeing up G12/H12 for the
# newly appended "related strategy" label
$personality = $ws.Range("H12").Value()
$ws.Range("H12").Value = ""
$ws.Range("G12").Value = "[Response Modulation, Lightly]"
$ws.Range("H13").Value = $personality
